$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 41599
$ws.Range("B14").NumberFormat = "m""月""d""日"""
$ws.Range("C14").Value = "雷建坤"
$ws.Range("E14").Value = "添加内容访问次数统计的功能"

$ws.Range("E15").Select() | Out-Null
